# Update column G ("K") values on Sheet1 per regenerated save_data.
# New K values replace the previous "Strike#" derived values, row by row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @(4, 4, 5, 9, 5, 4, 5, 6, 3, 5, 6, 6, 5, 9, 8, 6, 3, 7, 7, 7, 5, 4, 2, 5, 4)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Range("G$row").Value = $kValues[$i]
}
